$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <-> Row 3 swap (Fecha / Volumen)
$ws.Range("D2").Value = 44253
$ws.Range("M2").Value = 160

$ws.Range("D3").Value = 44250
$ws.Range("M3").Value = 200

# Row 4 <-> Row 5 swap (Fecha / Volumen / Precios)
$ws.Range("D4").Value = 44252
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 13500
$ws.Range("S4").Value = 750

$ws.Range("D5").Value = 44257
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("S5").Value = 806
